# Weekly data refresh: a new price record is inserted as row 64 (shifting
# all subsequent rows down by one) and the date in row 63 is advanced by
# one day to reflect the new latest entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing row 63 and insert it before row 64; this pushes rows
# 64..146 down to 65..147 and places a duplicate of the old row 63 into
# the newly created row 64.
$ws.Rows.Item(63).Copy() | Out-Null
$ws.Rows.Item(64).Insert() | Out-Null

# Row 63 keeps its original data except for the date, which moves one day
# forward (45174 -> 45175).
$ws.Range("D63").Value = 45175
